$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '77.354.67'
$ws.Range("E2").Value = '  +1.36%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.140.55'
$ws.Range("E3").Value = '  +5.69%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '202.68'
$ws.Range("E5").Value = '  +1.41%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.02'
$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.226'
$ws.Range("E8").Value = '  +13.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.569'
$ws.Range("E9").Value = '  +3.81%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.137.12'
$ws.Range("E10").Value = '  +5.57%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.530'
$ws.Range("E11").Value = '  +22.49%  '

$ws.Range("E12").Value = '  +1.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.39'
$ws.Range("E13").Value = '  +8.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.706.72'
$ws.Range("E14").Value = '  +5.35%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000224'
$ws.Range("E15").Value = '  +20.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.43'
$ws.Range("E16").Value = '  +4.93%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '77.180.90'
$ws.Range("E17").Value = '  +1.22%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.126.19'
$ws.Range("E18").Value = '  +5.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.02'
$ws.Range("E19").Value = '  +4.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.20'
$ws.Range("E20").Value = '  +5.23%  '

$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '431.04'
$ws.Range("E21").Value = '  +15.40%  '

$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.85'
$ws.Range("E22").Value = '  +27.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.71'
$ws.Range("E23").Value = '  +9.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.72'
$ws.Range("E24").Value = '  +5.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.302.74'
$ws.Range("E25").Value = '  +5.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.65'
$ws.Range("E26").Value = '  +8.35%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '75.25'
$ws.Range("E27").Value = '  +3.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.53'
$ws.Range("E28").Value = '  +9.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.14%  '

$ws.Range("E30").Value = '  +7.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.995'
$ws.Range("E31").Value = '  -0.17%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.73'
$ws.Range("E32").Value = '  +5.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.47'
$ws.Range("E33").Value = '  +5.72%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '524.04'
$ws.Range("E34").Value = '  +1.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.95'
$ws.Range("E35").Value = '  +0.88%  '

$ws.Range("E36").Value = '  +21.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '22.10'
$ws.Range("E37").Value = '  +9.11%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '163.74'
$ws.Range("E39").Value = '  +0.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.390'
$ws.Range("E40").Value = '  +1.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '194.97'
$ws.Range("E41").Value = '  +7.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.06'

$ws.Range("E43").Value = '  +0.03%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.106'
$ws.Range("E44").Value = '  +0.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.35'
$ws.Range("E45").Value = '  +8.64%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.794'
$ws.Range("E46").Value = '  +14.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.75'
$ws.Range("E47").Value = '  +7.36%  '

$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '42.61'
$ws.Range("E48").Value = '  -0.38%  '

$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("E49").Value = '  +4.82%  '

$ws.Range("E50").Value = '  +10.24%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.616'
$ws.Range("E51").Value = '  +5.58%  '
